# Simulated Wild Card round and logged it
#
# Rushing sheet: D.Johnson gets a new row inserted right after D.Mills
# (shifting the rest of the roster down one row), and J.Akins gains a
# fresh row at the bottom. Receiving sheet: stat totals bump for the
# players who touched the ball in the Wild Card game.

$wb = $excel.ActiveWorkbook

$rushing = $wb.Worksheets.Item("Rushing")
$receiving = $wb.Worksheets.Item("Receiving")

# ---- Rushing: rebuild rows 2..11 (rank, name, 1DATT, 2DATT, 3DATT, RZATT) ----
$names = @("D.Mills","D.Johnson","R.Burkhead","S.Phillips","R.Freeman","J.Samuels","T.Smith","C.Conley","P.Dorsett","J.Akins")
$stats = @(
    @(3,2,4,2),
    @(2,3,0,0),
    @(48,51,24,18),
    @(3,3,0,1),
    @(22,12,0,3),
    @(2,2,0,0),
    @(1,0,0,1),
    @(1,0,1,0),
    @(1,0,0,0),
    @(0,0,1,0)
)

for ($i = 0; $i -lt $names.Count; $i++) {
    $r = $i + 2
    $row = $stats[$i]
    $rushing.Cells.Item($r, 1).Value = $i
    $rushing.Cells.Item($r, 2).Value = $names[$i]
    $rushing.Cells.Item($r, 3).Value = $row[0]
    $rushing.Cells.Item($r, 4).Value = $row[1]
    $rushing.Cells.Item($r, 5).Value = $row[2]
    $rushing.Cells.Item($r, 6).Value = $row[3]
}

# New row 11 (J.Akins) needs the same bold/centered/bordered style as the
# rest of column A - copy the format down from the row above.
$rushing.Cells.Item(10, 1).Copy()
$rushing.Cells.Item(11, 1).PasteSpecial(-4122)

# ---- Receiving: bump stat totals for players who played the Wild Card game ----
$receiving.Cells.Item(2, 3).Value = 5
$receiving.Cells.Item(2, 4).Value = 5
$receiving.Cells.Item(2, 5).Value = 1

$receiving.Cells.Item(3, 3).Value = 45
$receiving.Cells.Item(3, 4).Value = 37
$receiving.Cells.Item(3, 5).Value = 3
$receiving.Cells.Item(3, 6).Value = 1

$receiving.Cells.Item(7, 3).Value = 104
$receiving.Cells.Item(7, 4).Value = 80
$receiving.Cells.Item(7, 5).Value = 32
$receiving.Cells.Item(7, 6).Value = 12

$receiving.Cells.Item(10, 3).Value = 29
$receiving.Cells.Item(10, 4).Value = 20
$receiving.Cells.Item(10, 5).Value = 8
$receiving.Cells.Item(10, 6).Value = 4
$receiving.Cells.Item(10, 7).Value = 3
$receiving.Cells.Item(10, 8).Value = 2

$receiving.Cells.Item(12, 3).Value = 40
$receiving.Cells.Item(12, 4).Value = 21
$receiving.Cells.Item(12, 5).Value = 16
$receiving.Cells.Item(12, 6).Value = 8

$receiving.Cells.Item(13, 3).Value = 7
$receiving.Cells.Item(13, 4).Value = 7
$receiving.Cells.Item(13, 5).Value = 1
$receiving.Cells.Item(13, 6).Value = 1

$receiving.Cells.Item(19, 3).Value = 4
$receiving.Cells.Item(19, 4).Value = 4

Write-Host "Wild Card round simulated and logged."
